$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the old first data row (A2/B2 = 44966 / 30.58) -----------------
# This shifts the remaining 17 historical rows up by one, so column A's
# dates line up with the new target layout (rows 2..18).
$ws.Rows("2:2").Delete()

# --- Header row ----------------------------------------------------------
# Order matters for shared-string table allocation (matches how the
# workbook's xl/sharedStrings.xml indices line up with the target file):
# "Ventas" -> 0, "Crecimiento" -> 1, "Fecha" -> 2, "Kpi" -> 3.
$ws.Range("B1").Value = "Ventas"
$ws.Range("C1").Value = "Crecimiento"
$ws.Range("A1").Value = "Fecha"
$ws.Range("D1").Value = "Kpi"

# --- New "Ventas" values for column B (rows 2..18) ----------------------
$ventas = @(513.98, 469.82, 386.06, 280.52, 232.89, 177.87, 135.99, 107.01, 88.99, 74.45, 61.09, 48.08, 34.2, 24.51, 19.17, 14.84, 10.71)
for ($i = 0; $i -lt $ventas.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $ventas[$i]
}

# --- Crecimiento (column C) and Kpi (column D, row 2 only) formulas -----
$ws.Range("C2:C18").NumberFormat = "0%"
$ws.Range("D2").NumberFormat = "0%"

$ws.Range("C2").Formula = "=(B2-B3)/100"
$ws.Range("C3:C18").Formula = "=(B3-B4)/100"
$ws.Range("D2").Formula = "=C2+C3/3"

# --- Selection / active cell --------------------------------------------
[void]$ws.Range("D2").Select()
